$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill column E (rows 2-9) with 1, mirroring the diff's added <c r="E2">..<c r="E9">
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Value = 1
}

# Update the active selection to E10, matching the diff's <selection activeCell="E10" sqref="E10"/>
$ws.Range("E10").Select()
